# Opencart_LoginData.xlsx: fix typo'd email domain and move the selection
# to the cell that was just edited (A2), matching the author's workflow of
# clicking into A2, retyping the address, then leaving the selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# kaw760@hotmail.com -> kaw760@gmail.com (cell A2, which carries a
# mailto: hyperlink whose displayed text is the cell value itself)
$ws.Range("A2").Value = "kaw760@gmail.com"

# Selection moves from C11 to A2
$ws.Range("A2").Select()
